$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# O12 was stored as a text "12"; change it to the numeric value 12.
$ws.Range("O12").Value = 12

# Append new row 13 with the latest survey submission.
$ws.Range("A13").Value = "pedro.paulistano@mrv.com.br"
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = "Painel do Portifólio - Planejamento da Produção - PLNESROBR004; PAP - Dossiê"
$ws.Range("G13").Value = ""
$ws.Range("H13").Value = ""
$ws.Range("I13").Value = "2025-05-20 13:02:13"
$ws.Range("J13").Value = "Painel do Portifólio - Planejamento da Produção - PLNESROBR004: asdfghjkl; PAP - Dossiê: qwertyui"
$ws.Range("K13").Value = "Planilha automatizada"
$ws.Range("L13").Value = "qqwwweee"
$ws.Range("M13").Value = "Painel Power BI"
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 7
$ws.Range("P13").Value = "Planilha geral - teste"
$ws.Range("Q13").Value = "wreqrweqarqrweqwe"
$ws.Range("R13").Value = "Painel Power BI"
$ws.Range("S13").Value = 4
$ws.Range("T13").Value = 4
